$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the scenario list (rows 17-18), pushing
# all the existing scenario rows (and the trailing blank rows) down by two.
$ws.Rows("17:18").Insert()

# Copy the formatting (yellow-highlighted "customFormat" style used by the
# two rows directly above, 15:16) onto the two freshly inserted rows so they
# match style s="4" instead of inheriting a generic copy of row 16's look.
$ws.Range("A16:D16").Copy()
$ws.Range("A17:D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New scenario row 17: Grant Type Profile (not executed)
$ws.Range("A17").Value = "MDOT"
$ws.Range("B17").Value = "TS_SM_FOU_GrantTypeProfile_Regression_001"
$ws.Range("C17").Value = "Grant Type Profile"
$ws.Range("D17").Value = "no"

# New scenario row 18: Grant Category Profile (executed)
$ws.Range("A18").Value = "MDOT"
$ws.Range("B18").Value = "TS_SM_FOU_GrantCategoryProfile_Regression_001"
$ws.Range("C18").Value = "Grant Category Profile"
$ws.Range("D18").Value = "yes"

# The previously-executed "NacuboGLAccountCategoryProfile" scenario row
# (shifted down to row 33 by the insert) is flipped back to not-executed.
$target = $ws.Range("C1:C43").Find("NacuboGLAccountCategoryProfile")
$target.Offset(0, 1).Value = "no"
